$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Locate the anchor paragraphs we need (by index, re-fetched fresh
# every time to avoid any stale-reference surprises).
#   - "PDF [TODO]" is an existing ListParagraph bullet item (numId=1)
#     we use as the template for the new bullet paragraphs' numbering.
#   - "Klant instellingen" is the Heading1 paragraph the new content
#     gets appended after.
# ------------------------------------------------------------------
$templateIdx = 0
$headingIdx = 0
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t.Contains("PDF [TODO]")) { $templateIdx = $i }
    if ($t.Contains("Klant instellingen")) { $headingIdx = $i }
}

$listTemplate = $d.Paragraphs.Item($templateIdx).Range.ListFormat.ListTemplate

# Remove the stray _GoBack bookmark currently sitting at the start of
# the "Klant instellingen" heading - it moves down to the new trailing
# empty bullet paragraph later on.
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()

# ------------------------------------------------------------------
# 1) New bullet paragraph "Done" right after the heading.
# ------------------------------------------------------------------
$d.Paragraphs.Item($headingIdx).Range.InsertParagraphAfter()
$doneIdx = $headingIdx + 1
$d.Paragraphs.Item($doneIdx).Style = "List Paragraph"
$d.Paragraphs.Item($doneIdx).Range.ListFormat.ApplyListTemplate($listTemplate, $true)
$d.Paragraphs.Item($doneIdx).Range.InsertAfter("Done")

# ------------------------------------------------------------------
# 2) New Heading1 paragraph "Order instellingen".
# ------------------------------------------------------------------
$d.Paragraphs.Item($doneIdx).Range.InsertParagraphAfter()
$orderIdx = $doneIdx + 1
$d.Paragraphs.Item($orderIdx).Style = "Heading 1"
$d.Paragraphs.Item($orderIdx).Range.InsertAfter("Order instellingen")

# ------------------------------------------------------------------
# 3) New bullet paragraph "Afrekenen uitgeschakel [TODO]".
# ------------------------------------------------------------------
$d.Paragraphs.Item($orderIdx).Range.InsertParagraphAfter()
$afrekenenIdx = $orderIdx + 1
$d.Paragraphs.Item($afrekenenIdx).Style = "List Paragraph"
$d.Paragraphs.Item($afrekenenIdx).Range.ListFormat.ApplyListTemplate($listTemplate, $true)
$d.Paragraphs.Item($afrekenenIdx).Range.InsertAfter("Afrekenen ")
$d.Paragraphs.Item($afrekenenIdx).Range.InsertAfter("uitgeschakel")
$d.Paragraphs.Item($afrekenenIdx).Range.InsertAfter(" [TODO]")

# ------------------------------------------------------------------
# 4) New bullet paragraph "Anoniem uitchecken toegestaan uit".
# ------------------------------------------------------------------
$d.Paragraphs.Item($afrekenenIdx).Range.InsertParagraphAfter()
$anoniemIdx = $afrekenenIdx + 1
$d.Paragraphs.Item($anoniemIdx).Style = "List Paragraph"
$d.Paragraphs.Item($anoniemIdx).Range.ListFormat.ApplyListTemplate($listTemplate, $true)
$d.Paragraphs.Item($anoniemIdx).Range.InsertAfter("Anoniem uitchecken toegestaa")
$d.Paragraphs.Item($anoniemIdx).Range.InsertAfter("n uit")

# ------------------------------------------------------------------
# 5) Convert the first trailing empty paragraph into a bullet
#    paragraph and re-home the _GoBack bookmark there (empty range).
# ------------------------------------------------------------------
$bookmarkHomeIdx = $anoniemIdx + 1
$d.Paragraphs.Item($bookmarkHomeIdx).Style = "List Paragraph"
$d.Paragraphs.Item($bookmarkHomeIdx).Range.ListFormat.ApplyListTemplate($listTemplate, $true)
$d.Bookmarks.Add("_GoBack", $d.Paragraphs.Item($bookmarkHomeIdx).Range)
